$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.856696666666667
$ws.Range("H2").Value = 8.57009
$ws.Range("I2").Value = 0.05747862151401942
$ws.Range("J2").Value = 0.05747862151401942
$ws.Range("O2").Value = 0.1205865197384776
$ws.Range("P2").Value = 0.1205865197384776
$ws.Range("Q2").Value = 0.4748267886822222
$ws.Range("R2").Value = 4.27344109814
$ws.Range("S2").Value = 0.006931146927740787
$ws.Range("T2").Value = 0.006931146927740788
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.856696666666667
$ws.Range("H3").Value = 8.57009
$ws.Range("I3").Value = 0.05747862151401942
$ws.Range("J3").Value = 0.05747862151401942
$ws.Range("M3").Value = 0.4346316666666667
$ws.Range("O3").Value = 0.3153182019998201
$ws.Range("P3").Value = 0.3153182019998201
$ws.Range("Q3").Value = 1.241610833394444
$ws.Range("R3").Value = 11.17449750055
$ws.Range("S3").Value = 0.01812405558922878
$ws.Range("T3").Value = 0.01812405558922878
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.856696666666667
$ws.Range("H4").Value = 8.57009
$ws.Range("I4").Value = 0.05747862151401942
$ws.Range("J4").Value = 0.05747862151401942
$ws.Range("M4").Value = 0.09482533333333333
$ws.Range("N4").Value = 0.284476
$ws.Range("O4").Value = 0.06879423637033719
$ws.Range("P4").Value = 0.0687942363703372
$ws.Range("Q4").Value = 0.2708872136488889
$ws.Range("R4").Value = 2.43798492284
$ws.Range("S4").Value = 0.0039541978746766
$ws.Range("T4").Value = 0.003954197874676601
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.856696666666667
$ws.Range("H5").Value = 8.57009
$ws.Range("I5").Value = 0.05747862151401942
$ws.Range("J5").Value = 0.05747862151401942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6827183333333333
$ws.Range("N5").Value = 2.048155
$ws.Range("O5").Value = 0.495301041891365
$ws.Range("P5").Value = 0.4953010418913651
$ws.Range("Q5").Value = 1.950319187105555
$ws.Range("R5").Value = 17.55287268395
$ws.Range("S5").Value = 0.02846922112237325
$ws.Range("T5").Value = 0.02846922112237325
$ws.Range("I6").Value = 0.7708435061432634
$ws.Range("J6").Value = 0.7708435061432632
$ws.Range("O6").Value = 0.1205865197384776
$ws.Range("P6").Value = 0.1205865197384776
$ws.Range("S6").Value = 0.09295333566882191
$ws.Range("T6").Value = 0.09295333566882191
$ws.Range("I7").Value = 0.7708435061432634
$ws.Range("J7").Value = 0.7708435061432632
$ws.Range("M7").Value = 0.4346316666666667
$ws.Range("O7").Value = 0.3153182019998201
$ws.Range("P7").Value = 0.3153182019998201
$ws.Range("S7").Value = 0.2430609883803311
$ws.Range("T7").Value = 0.243060988380331
$ws.Range("I8").Value = 0.7708435061432634
$ws.Range("J8").Value = 0.7708435061432632
$ws.Range("M8").Value = 0.09482533333333333
$ws.Range("N8").Value = 0.284476
$ws.Range("O8").Value = 0.06879423637033719
$ws.Range("P8").Value = 0.0687942363703372
$ws.Range("Q8").Value = 3.632857644081778
$ws.Range("R8").Value = 32.695718796736
$ws.Range("S8").Value = 0.05302959036615913
$ws.Range("T8").Value = 0.05302959036615913
$ws.Range("I9").Value = 0.7708435061432634
$ws.Range("J9").Value = 0.7708435061432632
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6827183333333333
$ws.Range("N9").Value = 2.048155
$ws.Range("O9").Value = 0.495301041891365
$ws.Range("P9").Value = 0.4953010418913651
$ws.Range("Q9").Value = 26.15565301823111
$ws.Range("R9").Value = 235.40087716408
$ws.Range("S9").Value = 0.3817995917279512
$ws.Range("T9").Value = 0.3817995917279512
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 2.081608
$ws.Range("H10").Value = 6.244823999999999
$ws.Range("I10").Value = 0.04188332620983732
$ws.Range("J10").Value = 0.04188332620983732
$ws.Range("O10").Value = 0.1205865197384776
$ws.Range("P10").Value = 0.1205865197384776
$ws.Range("Q10").Value = 0.3459951675893332
$ws.Range("R10").Value = 3.113956508304
$ws.Range("S10").Value = 0.005050564542715645
$ws.Range("T10").Value = 0.005050564542715645
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 2.081608
$ws.Range("H11").Value = 6.244823999999999
$ws.Range("I11").Value = 0.04188332620983732
$ws.Range("J11").Value = 0.04188332620983732
$ws.Range("M11").Value = 0.4346316666666667
$ws.Range("O11").Value = 0.3153182019998201
$ws.Range("P11").Value = 0.3153182019998201
$ws.Range("Q11").Value = 0.9047327543866666
$ws.Range("R11").Value = 8.14259478948
$ws.Range("S11").Value = 0.01320657511425784
$ws.Range("T11").Value = 0.01320657511425784
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 2.081608
$ws.Range("H12").Value = 6.244823999999999
$ws.Range("I12").Value = 0.04188332620983732
$ws.Range("J12").Value = 0.04188332620983732
$ws.Range("M12").Value = 0.09482533333333333
$ws.Range("N12").Value = 0.284476
$ws.Range("O12").Value = 0.06879423637033719
$ws.Range("P12").Value = 0.0687942363703372
$ws.Range("Q12").Value = 0.1973891724693333
$ws.Range("R12").Value = 1.776502552224
$ws.Range("S12").Value = 0.002881331443255487
$ws.Range("T12").Value = 0.002881331443255488
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 2.081608
$ws.Range("H13").Value = 6.244823999999999
$ws.Range("I13").Value = 0.04188332620983732
$ws.Range("J13").Value = 0.04188332620983732
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6827183333333333
$ws.Range("N13").Value = 2.048155
$ws.Range("O13").Value = 0.495301041891365
$ws.Range("P13").Value = 0.4953010418913651
$ws.Range("Q13").Value = 1.421151944413333
$ws.Range("R13").Value = 12.79036749972
$ws.Range("S13").Value = 0.02074485510960834
$ws.Range("T13").Value = 0.02074485510960834
$ws.Range("G14").Value = 4.573220666666667
$ws.Range("H14").Value = 13.719662
$ws.Range("I14").Value = 0.09201621679565497
$ws.Range("J14").Value = 0.09201621679565496
$ws.Range("O14").Value = 0.1205865197384776
$ws.Range("P14").Value = 0.1205865197384776
$ws.Range("Q14").Value = 0.7601393975168889
$ws.Range("R14").Value = 6.841254577651999
$ws.Range("S14").Value = 0.01109591534288928
$ws.Range("T14").Value = 0.01109591534288928
$ws.Range("G15").Value = 4.573220666666667
$ws.Range("H15").Value = 13.719662
$ws.Range("I15").Value = 0.09201621679565497
$ws.Range("J15").Value = 0.09201621679565496
$ws.Range("M15").Value = 0.4346316666666667
$ws.Range("O15").Value = 0.3153182019998201
$ws.Range("P15").Value = 0.3153182019998201
$ws.Range("Q15").Value = 1.987666520387778
$ws.Range("R15").Value = 17.88899868349
$ws.Range("S15").Value = 0.02901438803483157
$ws.Range("T15").Value = 0.02901438803483157
$ws.Range("G16").Value = 4.573220666666667
$ws.Range("H16").Value = 13.719662
$ws.Range("I16").Value = 0.09201621679565497
$ws.Range("J16").Value = 0.09201621679565496
$ws.Range("M16").Value = 0.09482533333333333
$ws.Range("N16").Value = 0.284476
$ws.Range("O16").Value = 0.06879423637033719
$ws.Range("P16").Value = 0.0687942363703372
$ws.Range("Q16").Value = 0.4336571741235555
$ws.Range("R16").Value = 3.902914567112
$ws.Range("S16").Value = 0.006330185368144479
$ws.Range("T16").Value = 0.006330185368144479
$ws.Range("G17").Value = 4.573220666666667
$ws.Range("H17").Value = 13.719662
$ws.Range("I17").Value = 0.09201621679565497
$ws.Range("J17").Value = 0.09201621679565496
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6827183333333333
$ws.Range("N17").Value = 2.048155
$ws.Range("O17").Value = 0.495301041891365
$ws.Range("P17").Value = 0.4953010418913651
$ws.Range("Q17").Value = 3.122221591512222
$ws.Range("R17").Value = 28.09999432361
$ws.Range("S17").Value = 0.04557572804978963
$ws.Range("T17").Value = 0.04557572804978963
$ws.Range("G18").Value = 1.877589
$ws.Range("H18").Value = 5.632767
$ws.Range("I18").Value = 0.037778329337225
$ws.Range("J18").Value = 0.037778329337225
$ws.Range("O18").Value = 0.1205865197384776
$ws.Range("P18").Value = 0.1205865197384776
$ws.Range("Q18").Value = 0.312084081498
$ws.Range("R18").Value = 2.808756733482
$ws.Range("S18").Value = 0.004555557256309991
$ws.Range("T18").Value = 0.004555557256309991
$ws.Range("G19").Value = 1.877589
$ws.Range("H19").Value = 5.632767
$ws.Range("I19").Value = 0.037778329337225
$ws.Range("J19").Value = 0.037778329337225
$ws.Range("M19").Value = 0.4346316666666667
$ws.Range("O19").Value = 0.3153182019998201
$ws.Range("P19").Value = 0.3153182019998201
$ws.Range("Q19").Value = 0.8160596363850001
$ws.Range("R19").Value = 7.344536727465001
$ws.Range("S19").Value = 0.01191219488117084
$ws.Range("T19").Value = 0.01191219488117084
$ws.Range("G20").Value = 1.877589
$ws.Range("H20").Value = 5.632767
$ws.Range("I20").Value = 0.037778329337225
$ws.Range("J20").Value = 0.037778329337225
$ws.Range("M20").Value = 0.09482533333333333
$ws.Range("N20").Value = 0.284476
$ws.Range("O20").Value = 0.06879423637033719
$ws.Range("P20").Value = 0.0687942363703372
$ws.Range("Q20").Value = 0.178043002788
$ws.Range("R20").Value = 1.602387025092
$ws.Range("S20").Value = 0.002598931318101501
$ws.Range("T20").Value = 0.002598931318101501
$ws.Range("G21").Value = 1.877589
$ws.Range("H21").Value = 5.632767
$ws.Range("I21").Value = 0.037778329337225
$ws.Range("J21").Value = 0.037778329337225
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.6827183333333333
$ws.Range("N21").Value = 2.048155
$ws.Range("O21").Value = 0.495301041891365
$ws.Range("P21").Value = 0.4953010418913651
$ws.Range("Q21").Value = 1.281864432765
$ws.Range("R21").Value = 11.536779894885
$ws.Range("S21").Value = 0.01871164588164267
$ws.Range("T21").Value = 0.01871164588164266
